# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-holdings layout, columns A-H) right
#    after "2021-Q4" and before "总计". It is created by copying the
#    existing "2021-Q4" sheet (same layout/styles) and then overwriting its
#    data, so the header row / column-A cells keep the workbook's existing
#    "bold + bordered + centered" style instead of minting new ones.
# 2) Insert a new first data row into "总计" for 2022-Q1 (date/count/value),
#    shifting the previous rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# The template sheet ("2021-Q4") only has 18 rows (1 header + 17 data); we
# need 22 (1 header + 21 data), so stretch column A's style down first.
$newSheet.Range("A18").Copy()
$newSheet.Range("A19:A22").PasteSpecial(-4122)

$fundRows = @(
    @(2, 0, "910007", "东方红启元三年持有期混合A", "69.33", "74.32", "4.60", "3.1892", 6),
    @(3, 1, "007887", "东方红启元三年持有期混合B", "59.71", "74.32", "4.60", "2.7467", 6),
    @(4, 2, "169107", "东方红恒阳五年定期开放混合", "22.26", "78.45", "4.39", "0.9772", 6),
    @(5, 3, "001302", "前海开源金银珠宝主题精选混合A", "8.61", "91.91", "8.70", "0.7491", 6),
    @(6, 4, "167508", "安信价值发现两年定期开放混合（LOF）", "5.43", "89.42", "5.67", "0.3079", 4),
    @(7, 5, "002207", "前海开源金银珠宝主题精选混合C", "3.45", "91.91", "8.70", "0.3002", 6),
    @(8, 6, "161609", "融通动力先锋混合", "7.12", "80.93", "3.21", "0.2286", 6),
    @(9, 7, "009766", "安信平稳双利3个月持有期混合A", "2.33", "39.45", "4.92", "0.1146", 3),
    @(10, 8, "001152", "融通新区域新经济灵活配置混合", "2.81", "80.98", "3.21", "0.0902", 5),
    @(11, 9, "510081", "长盛动态精选混合", "3.15", "60.76", "2.86", "0.0901", 10),
    @(12, 10, "003345", "安信新成长灵活配置混合A", "6.27", "30.76", "1.42", "0.0890", 6),
    @(13, 11, "001891", "中欧成长优选回报灵活配置混合E", "2.97", "94.42", "2.57", "0.0763", 10),
    @(14, 12, "166020", "中欧成长优选回报灵活配置混合A", "2.97", "94.42", "2.57", "0.0763", 10),
    @(15, 13, "001715", "工银瑞信新焦点灵活配置混合A", "1.42", "89.10", "4.86", "0.0690", 10),
    @(16, 14, "004393", "安信合作创新主题沪港深灵活配置混合", "0.49", "89.26", "6.08", "0.0298", 5),
    @(17, 15, "004249", "安信中国制造2025沪港深灵活配置混合", "0.58", "89.89", "4.50", "0.0261", 6),
    @(18, 16, "003346", "安信新成长灵活配置混合C", "1.18", "30.76", "1.42", "0.0168", 6),
    @(19, 17, "001998", "工银瑞信新焦点灵活配置混合C", "0.33", "89.10", "4.86", "0.0160", 10),
    @(20, 18, "009767", "安信平稳双利3个月持有期混合C", "0.26", "39.45", "4.92", "0.0128", 3),
    @(21, 19, "750005", "安信平稳增长混合A", "0.08", "65.16", "5.92", "0.0047", 2),
    @(22, 20, "002035", "安信平稳增长混合C", "0.00", "65.16", "5.92", "0", 2)
)

foreach ($r in $fundRows) {
    $rowIdx = $r[0]
    $newSheet.Cells.Item($rowIdx, 1).Value = $r[1]
    $newSheet.Cells.Item($rowIdx, 2).Value = "'" + $r[2]
    $newSheet.Cells.Item($rowIdx, 3).Value = $r[3]
    $newSheet.Cells.Item($rowIdx, 4).Value = "'" + $r[4]
    $newSheet.Cells.Item($rowIdx, 5).Value = "'" + $r[5]
    $newSheet.Cells.Item($rowIdx, 6).Value = "'" + $r[6]
    if ($rowIdx -eq 22) {
        $newSheet.Cells.Item($rowIdx, 7).Value = 0
    } else {
        $newSheet.Cells.Item($rowIdx, 7).Value = "'" + $r[7]
    }
    $newSheet.Cells.Item($rowIdx, 8).Value = $r[8]
}

# ---------------------------------------------------------------------
# Step 2: "总计" sheet gets a new first data row for 2022-Q1
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 21
$total.Range("D2").Value = 9.210000000000001

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
